$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "I"  = -0.1685440409879511
    "J"  = -0.1685440409879511
    "V"  = 1.670673076923077
    "X"  = 2.883122350314508
    "Y"  = -4.118049585241742
    "Z"  = 0.725156645441475
    "AA" = -0.1222208313719731
    "AB" = 0.08234438985867926
    "AC" = -0.2045652212306523
    "AE" = 307.5607606569383
    "AF" = 788.6607606569382
    "AG" = 774.7607606569383
    "AH" = 0.9895606011955144
    "AI" = 0.9998833601546784
    "AJ" = 0.9893752976474351
    "AK" = 0.9998812677650887
    "AP" = 10.45561080508689
}

foreach ($row in 2, 3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
